# Auto update Excel log
# Appends 6 new PRESENCE_DETECTED rows to the "mmWave" sheet (rows 68-73),
# continuing the existing log pattern (Living Room presence events).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# New rows to append: row number + timestamp (Date/Hour/Location/Value/Status
# are constant, following the pattern of the preceding rows).
$newRows = @(
    @{ Row = 68; Timestamp = "17:48:53" },
    @{ Row = 69; Timestamp = "17:49:01" },
    @{ Row = 70; Timestamp = "17:49:11" },
    @{ Row = 71; Timestamp = "17:49:22" },
    @{ Row = 72; Timestamp = "17:49:32" },
    @{ Row = 73; Timestamp = "17:49:43" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Force text storage for the date/time-looking values so Excel doesn't
    # auto-convert them into date/time serial numbers (they're plain text
    # in the rest of the log).
    $ws.Range("A" + $r).NumberFormat = "@"
    $ws.Range("A" + $r).Value = "2026-02-01"

    $ws.Range("B" + $r).NumberFormat = "@"
    $ws.Range("B" + $r).Value = $entry.Timestamp

    $ws.Range("C" + $r).NumberFormat = "@"
    $ws.Range("C" + $r).Value = "17:00"

    $ws.Range("D" + $r).Value = "Living Room"
    $ws.Range("E" + $r).Value = "PRESENCE_DETECTED"
    $ws.Range("F" + $r).Value = "Active"
}
